# "Generate Report for Archive"
# The handback/handoff status text moves from "Ready for handoff" to
# "In Translation" on every sheet that surfaces the localization status
# (the "Overview" roll-up sheet as well as each per-locale detail sheet).
# Excel's column-autofit then narrows the status column(s) to match the
# shorter replacement string.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# ColumnWidth (character units) snaps to Excel's internal pixel grid; 12.5
# is the value that lands the stored column width nearest the narrower
# autofit result for the new, shorter status text.
$newColumnWidth = 12.5

# --- Overview sheet: status mirrored in both the zh-cn and de-de columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E1").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").ColumnWidth = $newColumnWidth

# --- zh-cn detail sheet: Status column ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth

# --- de-de detail sheet: Status column ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
